# Financials update: insert a new "latest period" column before column D
# on the NEXA sheet, shifting the existing D:K data to E:L, then populate
# the new column D with the newest period's figures. Row 91 ("Capital
# Expenditures") also received restated values for its first five periods.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Insert a new blank column at D; existing D:K shift right to E:L.
$ws.Columns("D:D").Insert(-4121, 1)

# 2) Pick up the per-row number formatting (date header rows vs numeric
#    rows) by copying column E's formats (now holding what used to be in D)
#    into the freshly inserted column D.
$ws.Columns("E:E").Copy() | Out-Null
$ws.Columns("D:D").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# 3) Populate the new column D with the latest reporting period's values.
$newColD = @{
    7   = 43465
    8   = 2491200
    9   = 1888900
    10  = 602300
    12  = 126300
    14  = -34300
    17  = 2156600
    18  = 334600
    20  = -125000
    21  = 476700
    22  = 77600
    23  = 131900
    24  = 40900
    26  = 91000
    27  = 74900
    32  = 125000
    33  = 74900
    35  = 74900
    38  = 43465
    41  = 1032900
    42  = 91900
    43  = 173200
    44  = 269700
    45  = 130200
    46  = 1698000
    47  = 1400
    48  = 1968500
    49  = 1742500
    52  = 325100
    54  = 5735400
    57  = 387200
    58  = 32500
    59  = 232100
    60  = 651900
    61  = 1392400
    62  = 789400
    66  = 3258800
    72  = 1380200
    76  = 2476600
    80  = 43465
    81  = 74900
    83  = 267200
    89  = 347600
    94  = -158100
    96  = -3500
    100 = -177400
    101 = 1800
    102 = 13900
}

foreach ($row in $newColD.Keys) {
    $ws.Range("D$row").Value2 = $newColD[$row]
}

# Row 29 ("Discontinued Operations") has no data for the new period.
$ws.Range("D29").Value2 = "NA"

# 4) Row 91 ("Capital Expenditures") received restated figures across its
#    first five periods, not just a simple shift of the old values.
$ws.Range("D91").Value2 = -299800
$ws.Range("E91").Value2 = -196700
$ws.Range("F91").Value2 = -180900
$ws.Range("G91").Value2 = -183200
$ws.Range("H91").Value2 = -152000
